# Apply corrected IFRS figures to the NH투자증권 company_list sheet.
# (commit: "error solve ifrs list")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 48274
$ws.Range("E2").Value = 1255
$ws.Range("F2").Value = 1255
$ws.Range("G2").Value = 1199
$ws.Range("H2").Value = 813
$ws.Range("I2").Value = 811
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 410592
$ws.Range("L2").Value = 366811
$ws.Range("M2").Value = 43781
$ws.Range("N2").Value = 43689
$ws.Range("O2").Value = 91
$ws.Range("P2").Value = 15313
$ws.Range("Q2").Value = -8130
$ws.Range("R2").Value = -3383
$ws.Range("S2").Value = 15553
$ws.Range("T2").Value = 136
$ws.Range("V2").Value = 53051
$ws.Range("W2").Value = 2.6
$ws.Range("X2").Value = 1.68
$ws.Range("Y2").Value = 2.08
$ws.Range("Z2").Value = 0.23
$ws.Range("AA2").Value = 837.83
$ws.Range("AB2").Value = 187.25
$ws.Range("AC2").Value = 371
$ws.Range("AD2").Value = 27.6
$ws.Range("AE2").Value = 14596
$ws.Range("AF2").Value = 0.7
$ws.Range("AG2").Value = 160
$ws.Range("AH2").Value = 1.56
$ws.Range("AI2").Value = 60.23
$ws.Range("AJ2").Value = 281408887

# Row 3
$ws.Range("D3").Value = 70037
$ws.Range("E3").Value = 3141
$ws.Range("F3").Value = 3141
$ws.Range("G3").Value = 2822
$ws.Range("H3").Value = 2142
$ws.Range("I3").Value = 2151
$ws.Range("J3").Value = -8
$ws.Range("K3").Value = 417063
$ws.Range("L3").Value = 371558
$ws.Range("M3").Value = 45505
$ws.Range("N3").Value = 45423
$ws.Range("O3").Value = 82
$ws.Range("P3").Value = 15313
$ws.Range("Q3").Value = 6353
$ws.Range("R3").Value = -2001
$ws.Range("S3").Value = 1547
$ws.Range("T3").Value = 266
$ws.Range("V3").Value = 65347
$ws.Range("W3").Value = 4.49
$ws.Range("X3").Value = 3.06
$ws.Range("Y3").Value = 4.83
$ws.Range("Z3").Value = 0.52
$ws.Range("AA3").Value = 816.51
$ws.Range("AB3").Value = 198.49
$ws.Range("AC3").Value = 716
$ws.Range("AD3").Value = 14.1
$ws.Range("AE3").Value = 15175
$ws.Range("AF3").Value = 0.67
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 3.96
$ws.Range("AI3").Value = 56.11
$ws.Range("AJ3").Value = 281408887

# Row 4
$ws.Range("D4").Value = 88415
$ws.Range("E4").Value = 3019
$ws.Range("F4").Value = 3019
$ws.Range("G4").Value = 3020
$ws.Range("H4").Value = 2362
$ws.Range("I4").Value = 2361
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 429706
$ws.Range("L4").Value = 383218
$ws.Range("M4").Value = 46488
$ws.Range("N4").Value = 46388
$ws.Range("O4").Value = 100
$ws.Range("P4").Value = 15313
$ws.Range("Q4").Value = 22531
$ws.Range("R4").Value = -22956
$ws.Range("S4").Value = -5166
$ws.Range("T4").Value = 221
$ws.Range("V4").Value = 61834
$ws.Range("W4").Value = 3.42
$ws.Range("X4").Value = 2.67
$ws.Range("Y4").Value = 5.14
$ws.Range("Z4").Value = 0.5600000000000001
$ws.Range("AA4").Value = 824.34
$ws.Range("AB4").Value = 204.91
$ws.Range("AC4").Value = 786
$ws.Range("AD4").Value = 12.27
$ws.Range("AE4").Value = 15497
$ws.Range("AF4").Value = 0.62
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 4.15
$ws.Range("AI4").Value = 51.12
$ws.Range("AJ4").Value = 281408887

# Row 5
$ws.Range("D5").Value = 95455
$ws.Range("E5").Value = 4592
$ws.Range("F5").Value = 4592
$ws.Range("G5").Value = 4425
$ws.Range("H5").Value = 3496
$ws.Range("I5").Value = 3501
$ws.Range("J5").Value = -5
$ws.Range("K5").Value = 438927
$ws.Range("L5").Value = 390565
$ws.Range("M5").Value = 48362
$ws.Range("N5").Value = 48306
$ws.Range("O5").Value = 56
$ws.Range("P5").Value = 15313
$ws.Range("Q5").Value = 8763
$ws.Range("R5").Value = -21733
$ws.Range("S5").Value = 12210
$ws.Range("T5").Value = 152
$ws.Range("V5").Value = 68219
$ws.Range("W5").Value = 4.81
$ws.Range("X5").Value = 3.66
$ws.Range("Y5").Value = 7.4
$ws.Range("Z5").Value = 0.8100000000000001
$ws.Range("AA5").Value = 807.59
$ws.Range("AB5").Value = 217.15
$ws.Range("AC5").Value = 1166
$ws.Range("AD5").Value = 11.92
$ws.Range("AE5").Value = 16138
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 3.6
$ws.Range("AI5").Value = 43.02
$ws.Range("AJ5").Value = 281408887

# Row 6
$ws.Range("D6").Value = 92413
$ws.Range("E6").Value = 5401
$ws.Range("F6").Value = 5401
$ws.Range("G6").Value = 5048
$ws.Range("H6").Value = 3615
$ws.Range("I6").Value = 3609
$ws.Range("K6").Value = 534226
$ws.Range("L6").Value = 483688
$ws.Range("M6").Value = 50538
$ws.Range("N6").Value = 50476
$ws.Range("P6").Value = 15313
$ws.Range("Q6").Value = -15955
$ws.Range("R6").Value = 9603
$ws.Range("S6").Value = 7293
$ws.Range("T6").Value = 165
$ws.Range("V6").Value = 81367
$ws.Range("W6").Value = 5.85
$ws.Range("X6").Value = 3.91
$ws.Range("Y6").Value = 7.31
$ws.Range("Z6").Value = 0.74
$ws.Range("AA6").Value = 957.08
$ws.Range("AB6").Value = 231.36
$ws.Range("AC6").Value = 1202
$ws.Range("AD6").Value = 10.86
$ws.Range("AE6").Value = 16863
$ws.Range("AF6").Value = 0.77
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 3.83
$ws.Range("AI6").Value = 41.74
$ws.Range("AJ6").Value = 281408887

# Row 7
$ws.Range("D7").Value = 134038
$ws.Range("E7").Value = 6335
$ws.Range("G7").Value = 6624
$ws.Range("H7").Value = 4937
$ws.Range("I7").Value = 4970
$ws.Range("K7").Value = 556784
$ws.Range("L7").Value = 501773
$ws.Range("M7").Value = 54225
$ws.Range("N7").Value = 53965
$ws.Range("P7").Value = 15311
$ws.Range("W7").Value = 4.73
$ws.Range("X7").Value = 3.68
$ws.Range("Y7").Value = 9.52
$ws.Range("Z7").Value = 0.91
$ws.Range("AA7").Value = 925.35
$ws.Range("AC7").Value = 1655
$ws.Range("AD7").Value = 7.22
$ws.Range("AE7").Value = 18013
$ws.Range("AF7").Value = 0.66
$ws.Range("AG7").Value = 621
$ws.Range("AH7").Value = 5.19
$ws.Range("AI7").Value = 35.15

# Row 8
$ws.Range("D8").Value = 103900
$ws.Range("E8").Value = 6163
$ws.Range("G8").Value = 5916
$ws.Range("H8").Value = 4378
$ws.Range("I8").Value = 4397
$ws.Range("K8").Value = 593826
$ws.Range("L8").Value = 536122
$ws.Range("M8").Value = 57090
$ws.Range("N8").Value = 56703
$ws.Range("P8").Value = 15311
$ws.Range("W8").Value = 5.93
$ws.Range("X8").Value = 4.21
$ws.Range("Y8").Value = 7.96
$ws.Range("Z8").Value = 0.76
$ws.Range("AA8").Value = 939.08
$ws.Range("AC8").Value = 1464
$ws.Range("AD8").Value = 7.55
$ws.Range("AE8").Value = 18927
$ws.Range("AF8").Value = 0.58
$ws.Range("AG8").Value = 635
$ws.Range("AH8").Value = 5.75
$ws.Range("AI8").Value = 40.65

# Row 9
$ws.Range("D9").Value = 124350
$ws.Range("E9").Value = 6619
$ws.Range("G9").Value = 6465
$ws.Range("H9").Value = 4745
$ws.Range("I9").Value = 4691
$ws.Range("K9").Value = 633133
$ws.Range("L9").Value = 572904
$ws.Range("M9").Value = 60230
$ws.Range("N9").Value = 60476
$ws.Range("P9").Value = 15311
$ws.Range("W9").Value = 5.32
$ws.Range("X9").Value = 3.82
$ws.Range("Y9").Value = 8.01
$ws.Range("Z9").Value = 0.77
$ws.Range("AA9").Value = 951.2
$ws.Range("AC9").Value = 1562
$ws.Range("AD9").Value = 7.07
$ws.Range("AE9").Value = 20187
$ws.Range("AF9").Value = 0.55
$ws.Range("AG9").Value = 655
$ws.Range("AH9").Value = 5.93
$ws.Range("AI9").Value = 39.32

# These columns were removed / no longer populated for this row in the corrected data
$ws.Range("U2").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("U6").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("U9").ClearContents()
